$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 100.125
$ws.Cells.Item(8, 9).Value = 107.28571
$ws.Cells.Item(8, 10).Value = 50
$ws.Cells.Item(8, 11).Value = 321.85713
$ws.Cells.Item(8, 12).Value = 150
$ws.Cells.Item(8, 13).Value = -182.85713
$ws.Cells.Item(8, 14).Value = -428

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 27031194
$ws.Cells.Item(132, 9).Value = 37038050
$ws.Cells.Item(132, 10).Value = 12690.2
$ws.Cells.Item(132, 11).Value = 111114150
$ws.Cells.Item(132, 12).Value = 38070.60000000001
$ws.Cells.Item(132, 13).Value = -111111620
$ws.Cells.Item(132, 14).Value = -43130.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2419091.8
$ws.Cells.Item(134, 9).Value = 936.12
$ws.Cells.Item(134, 10).Value = 5297848
$ws.Cells.Item(134, 11).Value = 2808.36
$ws.Cells.Item(134, 12).Value = 15893544
$ws.Cells.Item(134, 13).Value = -273.3600000000001
$ws.Cells.Item(134, 14).Value = -15898614

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 9528112
$ws.Cells.Item(132, 9).Value = 1232.762
$ws.Cells.Item(132, 11).Value = 3698.286
$ws.Cells.Item(132, 13).Value = -1168.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 130.53847
$ws.Cells.Item(2, 9).Value = 31.25
$ws.Cells.Item(2, 10).Value = 289.4
$ws.Cells.Item(2, 11).Value = 187.5
$ws.Cells.Item(2, 12).Value = 1736.4
$ws.Cells.Item(2, 13).Value = -74.5
$ws.Cells.Item(2, 14).Value = -1962.4

$ws.Cells.Item(9, 8).Value = 157171980
$ws.Cells.Item(9, 10).Value = 200040720
$ws.Cells.Item(9, 12).Value = 600122160
$ws.Cells.Item(9, 14).Value = -600122608

$ws.Cells.Item(15, 8).Value = 375
$ws.Cells.Item(15, 9).Value = 250
$ws.Cells.Item(15, 10).Value = 500
$ws.Cells.Item(15, 11).Value = 750
$ws.Cells.Item(15, 12).Value = 1500
$ws.Cells.Item(15, 13).Value = -610
$ws.Cells.Item(15, 14).Value = -1780

$ws.Cells.Item(20, 8).Value = 1000800

$ws.Cells.Item(21, 8).Value = 475.75
$ws.Cells.Item(21, 9).Value = 300.33334
$ws.Cells.Item(21, 10).Value = 1002
$ws.Cells.Item(21, 11).Value = 901.0000200000001
$ws.Cells.Item(21, 12).Value = 3006
$ws.Cells.Item(21, 13).Value = -728.0000200000001
$ws.Cells.Item(21, 14).Value = -3352

$ws.Cells.Item(22, 8).Value = 25038376
$ws.Cells.Item(22, 10).Value = 76000
$ws.Cells.Item(22, 12).Value = 228000
$ws.Cells.Item(22, 14).Value = -228338

$ws.Cells.Item(27, 8).Value = 25038376
$ws.Cells.Item(27, 10).Value = 76000
$ws.Cells.Item(27, 12).Value = 228000
$ws.Cells.Item(27, 14).Value = -228204

$ws.Cells.Item(34, 8).Value = 233.33333
$ws.Cells.Item(34, 9).Value = 233.33333
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 699.99999
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -615.99999
$ws.Cells.Item(34, 14).ClearContents()

$ws.Cells.Item(40, 8).Value = 114.21429
$ws.Cells.Item(40, 9).Value = 114.21429
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 456.85716
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -387.85716
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 50495
$ws.Cells.Item(46, 9).Value = 100000
$ws.Cells.Item(46, 10).Value = 990
$ws.Cells.Item(46, 11).Value = 300000
$ws.Cells.Item(46, 12).Value = 2970
$ws.Cells.Item(46, 13).Value = -299909
$ws.Cells.Item(46, 14).Value = -3152

$ws.Cells.Item(50, 8).Value = 1809.8572
$ws.Cells.Item(50, 9).Value = 42.25
$ws.Cells.Item(50, 10).Value = 4166.6665
$ws.Cells.Item(50, 11).Value = 126.75
$ws.Cells.Item(50, 12).Value = 12499.9995
$ws.Cells.Item(50, 13).Value = 354.25
$ws.Cells.Item(50, 14).Value = -13461.9995

$ws.Cells.Item(53, 8).Value = 1809.8572
$ws.Cells.Item(53, 9).Value = 42.25
$ws.Cells.Item(53, 10).Value = 4166.6665
$ws.Cells.Item(53, 11).Value = 126.75
$ws.Cells.Item(53, 12).Value = 12499.9995
$ws.Cells.Item(53, 13).Value = 354.25
$ws.Cells.Item(53, 14).Value = -13461.9995

$ws.Cells.Item(58, 8).Value = 28929.166
$ws.Cells.Item(58, 10).Value = 31527.273
$ws.Cells.Item(58, 12).Value = 94581.819
$ws.Cells.Item(58, 14).Value = -94837.819

$ws.Cells.Item(64, 8).Value = 2004.6666
$ws.Cells.Item(64, 9).Value = 1000
$ws.Cells.Item(64, 11).Value = 3000
$ws.Cells.Item(64, 13).Value = -2730

$ws.Cells.Item(67, 8).Value = 2004.6666
$ws.Cells.Item(67, 9).Value = 1000
$ws.Cells.Item(67, 11).Value = 3000
$ws.Cells.Item(67, 13).Value = -2064

$ws.Cells.Item(76, 8).Value = 2000
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 13).ClearContents()

$ws.Cells.Item(79, 8).Value = 2000
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 13).ClearContents()

$ws.Cells.Item(81, 8).Value = 778.25
$ws.Cells.Item(81, 9).Value = 778.25
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 2334.75
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -1211.75
$ws.Cells.Item(81, 14).ClearContents()

$ws.Cells.Item(82, 8).Value = 1400
$ws.Cells.Item(82, 10).Value = 2000
$ws.Cells.Item(82, 12).Value = 6000
$ws.Cells.Item(82, 14).Value = -6812

$ws.Cells.Item(84, 8).Value = 778.25
$ws.Cells.Item(84, 9).Value = 778.25
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 7004.25
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -1388.25
$ws.Cells.Item(84, 14).ClearContents()

$ws.Cells.Item(85, 8).Value = 1400
$ws.Cells.Item(85, 10).Value = 2000
$ws.Cells.Item(85, 12).Value = 6000
$ws.Cells.Item(85, 14).Value = -8808

$ws.Cells.Item(86, 8).Value = 142.85715
$ws.Cells.Item(86, 9).Value = 120
$ws.Cells.Item(86, 10).Value = 200
$ws.Cells.Item(86, 11).Value = 360
$ws.Cells.Item(86, 12).Value = 600
$ws.Cells.Item(86, 13).Value = 826
$ws.Cells.Item(86, 14).Value = -2972

$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).ClearContents()
$ws.Cells.Item(87, 14).ClearContents()

$ws.Cells.Item(88, 8).Value = 3000
$ws.Cells.Item(88, 10).Value = 3000
$ws.Cells.Item(88, 12).Value = 9000
$ws.Cells.Item(88, 14).Value = -9856

$ws.Cells.Item(89, 8).Value = 142.85715
$ws.Cells.Item(89, 9).Value = 120
$ws.Cells.Item(89, 10).Value = 200
$ws.Cells.Item(89, 11).Value = 1080
$ws.Cells.Item(89, 12).Value = 1800
$ws.Cells.Item(89, 13).Value = 4848
$ws.Cells.Item(89, 14).Value = -13656

$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).ClearContents()
$ws.Cells.Item(90, 14).ClearContents()

$ws.Cells.Item(91, 8).Value = 3000
$ws.Cells.Item(91, 10).Value = 3000
$ws.Cells.Item(91, 12).Value = 9000
$ws.Cells.Item(91, 14).Value = -11964

$ws.Cells.Item(92, 8).Value = 50138.5
$ws.Cells.Item(92, 9).Value = 277
$ws.Cells.Item(92, 10).Value = 100000
$ws.Cells.Item(92, 11).Value = 831
$ws.Cells.Item(92, 12).Value = 300000
$ws.Cells.Item(92, 13).Value = 417
$ws.Cells.Item(92, 14).Value = -302496

$ws.Cells.Item(94, 8).Value = 2000
$ws.Cells.Item(94, 10).Value = 2000
$ws.Cells.Item(94, 12).Value = 6000
$ws.Cells.Item(94, 14).Value = -7352

$ws.Cells.Item(102, 8).Value = 3200
$ws.Cells.Item(102, 9).Value = 3200
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 9600
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -7166
$ws.Cells.Item(102, 14).ClearContents()

$ws.Cells.Item(103, 8).Value = 4621.875
$ws.Cells.Item(103, 9).Value = 95
$ws.Cells.Item(103, 10).Value = 7338
$ws.Cells.Item(103, 11).Value = 285
$ws.Cells.Item(103, 12).Value = 22014
$ws.Cells.Item(103, 13).Value = 594
$ws.Cells.Item(103, 14).Value = -23772

$ws.Cells.Item(106, 8).Value = 26500
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 26500
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 79500
$ws.Cells.Item(106, 14).Value = -81392
$ws.Cells.Item(106, 13).ClearContents()

$ws.Cells.Item(109, 8).Value = 4385.4546
$ws.Cells.Item(109, 9).Value = 811.75
$ws.Cells.Item(109, 10).Value = 6427.5713
$ws.Cells.Item(109, 11).Value = 2435.25
$ws.Cells.Item(109, 12).Value = 19282.7139
$ws.Cells.Item(109, 13).Value = -1395.25
$ws.Cells.Item(109, 14).Value = -21362.7139

$ws.Cells.Item(113, 8).Value = 4365600.5
$ws.Cells.Item(113, 10).Value = 2632110.5
$ws.Cells.Item(113, 12).Value = 7896331.5
$ws.Cells.Item(113, 14).Value = -7900671.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2791.4443
$ws.Cells.Item(126, 9).Value = 3041.3333
$ws.Cells.Item(126, 10).Value = 2666.5
$ws.Cells.Item(126, 11).Value = 9123.999899999999
$ws.Cells.Item(126, 12).Value = 7999.5
$ws.Cells.Item(126, 13).Value = -6653.999899999999
$ws.Cells.Item(126, 14).Value = -12939.5

$ws.Cells.Item(132, 8).Value = 4401.8857
$ws.Cells.Item(132, 9).Value = 1329.6957
$ws.Cells.Item(132, 10).Value = 10290.25
$ws.Cells.Item(132, 11).Value = 3989.0871
$ws.Cells.Item(132, 12).Value = 30870.75
$ws.Cells.Item(132, 13).Value = -1459.0871
$ws.Cells.Item(132, 14).Value = -35930.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 23263474
$ws.Cells.Item(132, 9).Value = 62501532
$ws.Cells.Item(132, 10).Value = 11292.889
$ws.Cells.Item(132, 11).Value = 187504596
$ws.Cells.Item(132, 12).Value = 33878.667
$ws.Cells.Item(132, 13).Value = -187502066
$ws.Cells.Item(132, 14).Value = -38938.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 80005
$ws.Cells.Item(11, 10).Value = 80005
$ws.Cells.Item(11, 12).Value = 80005
$ws.Cells.Item(11, 14).Value = -80289

$ws.Cells.Item(132, 8).Value = 14698.329
$ws.Cells.Item(132, 9).Value = 18592.793
$ws.Cells.Item(132, 11).Value = 55778.379
$ws.Cells.Item(132, 13).Value = -53248.379
